$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: (row, newK, newL) - newK is the newly computed lower-bound viral load value,
# newL is the previous K value (now re-labeled as the upper-bound viral load column).
$rowData = @(
    @(3, "2.8327548408945739e-13", "8.5659226251169461e-13"),
    @(4, "3.1452590056913757e-10", "3600.567088612313"),
    @(5, "3.0112462200412679e-13", "2.7817775361036009e-13"),
    @(6, "3.0680631149652922e-13", "4.0762454291123869e-13"),
    @(7, "2.6564295986491696e-13", "3.5380707520490049e-13"),
    @(8, "2.1003619306277225e-13", "1.600446329362164e-13"),
    @(9, "1.3613161652797686e-13", "1.6046663974841837e-13"),
    @(10, "2.4871032810403427e-13", "3.1944496034882083e-13"),
    @(11, "2.5521551668671937e-13", "3.1460200513970353e-13"),
    @(12, "2.4538169275658034e-10", "434.79630932244839"),
    @(13, "8.9594554769161624e-11", "434.79630934467366"),
    @(14, "8.8051477701597343e-13", "4.3181780431108806e-14"),
    @(15, "5.0845433700672674e-11", "3600.5670885991171"),
    @(16, "1.8238573135827519e-13", "1.6358643185609223e-13"),
    @(17, "1459166.7937527974", "18354129.898997303"),
    @(18, "2.009198600874557e-10", "4025.2436895848841"),
    @(19, "1.4686981349029317e-13", "1.9253407447799867e-13"),
    @(20, "5.6316125705162113e-13", "4.8825797274677599e-13"),
    @(21, "2.6713670404378322e-13", "3.6593186067734761e-13"),
    @(22, "1.8709254947561437e-13", "1.5567031766016315e-13"),
    @(23, "2.8473100029579008e-13", "4.2173545850027153e-13"),
    @(24, "1.7659778084492708e-13", "1.5296503445595946e-13"),
    @(25, "1.6197156807536824e-13", "5.4991907059751861e-13"),
    @(26, "1.3442589391876798e-13", "1.7004285820233747e-13"),
    @(27, "2.7334311446556845e-13", "1.6981786278665822e-14"),
    @(28, "7.7574408014227799e-11", "434.79630935163527"),
    @(29, "1.0151687079091582e-12", "1.7760365182000824e-13"),
    @(30, "1.0213921225217977e-10", "4025.2436895615451"),
    @(31, "3.1343998361363802e-13", "7.3208421125555145e-13"),
    @(32, "1.7338736606294788e-10", "4025.2436895679998"),
    @(33, "1.5406410159382618e-13", "1.9224909126675068e-13"),
    @(34, "6.4776858930525085e-12", "4025.2436895197206"),
    @(35, "5.1263037761038787e-11", "3600.5670886465159"),
    @(36, "5.2388275286375982e-13", "1.0682448559336435e-13"),
    @(37, "9.6024214399731268e-11", "434.79630936436013"),
    @(38, "3.3500986843311463e-13", "5.4299675064695717e-13"),
    @(39, "1.9459171445763971e-13", "1.6235677398267961e-13"),
    @(40, "1.3594806183767444e-13", "1.521807854115385e-13"),
    @(41, "3.6270378390810211e-13", "3.7209480512525169e-13"),
    @(42, "11917.743108492341", "178213.45764733746"),
    @(45, "3.5019363930215738e-13", "3.2617156212824624e-13"),
    @(46, "2.9478039162414429e-10", "15247.710502793434"),
    @(47, "6.1779973855586275e-10", "9518.4154950389511"),
    @(48, "3.1647462224576241e-13", "3.4348170733076685e-13"),
    @(49, "2.8182266120580528e-13", "2.7123394351086998e-13"),
    @(50, "1.379483553323944e-10", "9518.415494980698"),
    @(51, "4.1938301016853237e-13", "3.6083562881035923e-13"),
    @(52, "1.0375692737371355e-12", "9.3209166390504672e-13"),
    @(53, "2.6903325911729091e-13", "1.1428668699964137e-13"),
    @(54, "2.3795415709585073e-13", "2.4363859544298442e-13"),
    @(55, "4.5332053122202978e-13", "3.8702888279161927e-13"),
    @(56, "7.8543744410075232e-13", "6.7293331304656518e-13"),
    @(57, "91318.775854704261", "1963985.5542124293"),
    @(58, "4.1473067030742193e-13", "4.2715901020859238e-13"),
    @(59, "3.2577207446795497e-13", "3.6435476690809524e-13"),
    @(60, "4.9351011089388793e-13", "4.8962242380909324e-13"),
    @(61, "2.5613235117319269e-10", "131885.34444839912"),
    @(62, "4.4203785452772411e-13", "3.5918466077560949e-13"),
    @(63, "9.4928553682028103e-13", "9.5848384556065047e-13"),
    @(64, "6.8572134626674034e-10", "9518.4154950405318"),
    @(65, "2.7472035110669624e-13", "2.6254334672668488e-13"),
    @(66, "1.1699979613332134e-12", "1.1758032473659328e-12"),
    @(67, "4.254274549025055e-13", "3.2528205449668707e-13"),
    @(68, "9.0827253191488133e-13", "1.0171369666762214e-12"),
    @(69, "5.0479613038586527e-13", "4.548785605567772e-13"),
    @(70, "7.43211646384424e-13", "4.1775351724258947e-13"),
    @(71, "3.4994360556824244e-10", "131885.34444890189"),
    @(72, "4.8105664322093626e-13", "6.3007379942264341e-13"),
    @(73, "7.5887747949432415e-13", "1.0301463737060455e-12"),
    @(74, "4.7438446891667184e-13", "5.1937025443116633e-13"),
    @(75, "7.8193566309065878e-13", "6.8554606213501644e-13"),
    @(76, "4.9923516482603172e-13", "6.5774166128526359e-13"),
    @(77, "5.6424998261200289e-13", "1.523663577204127e-13"),
    @(78, "6.6115891806186284e-13", "7.9579874449715471e-13"),
    @(79, "4.5570241610873168e-13", "3.4932844677401742e-13"),
    @(80, "86192.115961475502", "2320454.0659286371"),
    @(81, "4.3069134752849785e-13", "3.3737526096565995e-13"),
    @(82, "4.9040682668934789e-13", "4.0859190230596777e-13"),
    @(83, "2.9021491514401724e-13", "3.0853985550473767e-13"),
    @(84, "2.8894141923458364e-13", "2.9051085625237939e-13")
)

foreach ($item in $rowData) {
    $row = $item[0]
    $kVal = [double]$item[1]
    $lVal = [double]$item[2]
    $ws.Cells.Item($row, 11).Value = $kVal
    $ws.Cells.Item($row, 12).Value = $lVal
}

# Column L (12) width changes from 17.3984375 to 11.625
$ws.Columns.Item(12).ColumnWidth = 11.625

Write-Host "Applied COVID-19 decoding export fix to $($rowData.Count) rows"
